$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reset "Points Today" column (B2:B14) back to 0 for all students
$ws.Range("B2:B14").Value = 0

# Update the active cell selection to match the saved view state
$ws.Range("D16").Select()
